$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) keeps bold style s=1; just update text content
$ws.Range("A1").Value = ' **Balance Sheet Indicators**'
$ws.Range("B1").Value = '**Current Year (2024)**'
$ws.Range("C1").Value = '**Previous Year (2023)** '

# Force data rows 2-17 to text type so numeric-looking strings are not auto-converted
$ws.Range("A2:C17").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = ' Cash and Cash Equivalents'
$ws.Range("B2").Value = '1000'
$ws.Range("C2").Value = '1000                     '

# Row 3
$ws.Range("A3").Value = ' Accounts Receivable'
$ws.Range("B3").Value = '11,987,605.97'
$ws.Range("C3").Value = '10,711,454.12            '

# Row 4
$ws.Range("A4").Value = ' Property, Plant and Equipment (Net)'
$ws.Range("B4").Value = '3,494,523.92'
$ws.Range("C4").Value = '3,494,523.92             '

# Row 5
$ws.Range("A5").Value = ' Total Assets'
$ws.Range("B5").Value = '14,355,193.96'
$ws.Range("C5").Value = '13,424,369.47            '

# Row 6
$ws.Range("A6").Value = ' Accounts Payable'
$ws.Range("B6").Value = '-12,443,892.15'
$ws.Range("C6").Value = '-10,979,515.78           '

# Row 7
$ws.Range("A7").Value = ' Accumulated Profit/(Loss)'
$ws.Range("B7").Value = '-2,444,853.69'
$ws.Range("C7").Value = '-2,741,596.38            '

# Row 8
$ws.Range("A8").Value = ' Total Liabilities'
$ws.Range("B8").Value = '-14,888,745.84'
$ws.Range("C8").Value = '-13,721,112.16           '

# Row 9
$ws.Range("A9").Value = ' **Income Statement Indicators**'
$ws.Range("B9").Value = '**Current Year (2024)**'
$ws.Range("C9").Value = '**Previous Year (2023)** '

# Row 10
$ws.Range("A10").Value = ' Revenue'
$ws.Range("B10").Value = '-1,276,151.85'
$ws.Range("C10").Value = '-1,727,145.61            '

# Row 11
$ws.Range("A11").Value = ' Cost of Goods Sold'
$ws.Range("B11").Value = '367,148.33'
$ws.Range("C11").Value = '428,073.69               '

# Row 12
$ws.Range("A12").Value = ' Gross Profit'
$ws.Range("B12").Value = '*N/A*'
$ws.Range("C12").Value = '*N/A*                    '

# Row 13
$ws.Range("A13").Value = ' General and Administrative Expenses'
$ws.Range("B13").Value = '*Consolidated within Expense Total*'
$ws.Range("C13").Value = '*Consolidated within Expense Total* '

# Row 14
$ws.Range("A14").Value = ' Net Profit'
$ws.Range("B14").Value = '*N/A*'
$ws.Range("C14").Value = '*N/A*                    '

# Row 15
$ws.Range("A15").Value = ' **Cash Flow Indicators**'
$ws.Range("B15").Value = '**Current Year (2024)**'
$ws.Range("C15").Value = '**Previous Year (2023)** '

# Row 16
$ws.Range("A16").Value = ' Depreciation and Amortization'
$ws.Range("B16").Value = '350,277.36'
$ws.Range("C16").Value = '349,452.36               '

# Row 17
$ws.Range("A17").Value = ' Total Cash Flow'
$ws.Range("B17").Value = '533,551.88'
$ws.Range("C17").Value = '296,742.69               '

# Reset style of data rows back to Normal (removes temporary text format) while keeping text type
$ws.Range("A2:C17").Style = "Normal"

Write-Host "Applied financial summary table update"
